$d = $word.ActiveDocument

# --- 1. Insert the evaluator's name as a bold run after "Nombre y Apellidos: " ---
$target = $d.Paragraphs(40)
$insertAt = $target.Range.End - 1
$newRun = $d.Range($insertAt, $insertAt)
$newRun.InsertAfter("Ana Blasco Parra")
$fmtRange = $d.Range($insertAt, $insertAt + 17)
$fmtRange.Font.Bold = $true
$fmtRange.Font.Size = 9

# --- 2. Mint the footnotes/endnotes parts (with their separator scaffolding) ---
#        the same way Word does the first time a footnote is touched, then
#        remove the placeholder footnote itself so only the standard
#        separator/continuationSeparator entries remain.
$fnAnchor = $d.Paragraphs(1).Range
$tempNote = $d.Footnotes.Add($fnAnchor, "", "x")
$tempNote.Delete()

Write-Host "done"
